# "neue secure testversion hinzugefügt"
#
# The "Insecure IT Testmodell" tag-matrix sheet gets a secure-variant
# counterpart: two firewall hops (firewall1 / firewall2, reached via
# to-firewall1 / to-firewall2) are spliced into the existing
# app/client/database/webapp chain. That turns the original 10-row table
# into 14 rows while keeping the existing header + per-column styling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank rows at the positions where the new tags land, so the
# untouched rows above/below keep their original formatting and the
# "user-data / X" row ends up at row 13 (was row 9), "dmz" at row 14 (was
# row 10) - exactly matching the target layout.
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(9).Insert()
$ws.Rows.Item(11).Insert()
$ws.Rows.Item(13).Insert()

# Full resulting A:B table (column A = tag/element name, column B = flag).
$values = @(
    @("Element", "pii"),
    @("app", $null),
    @("to-database", $null),
    @("client", $null),
    @("to-firewall1", $null),
    @("database", $null),
    @("firewall1", $null),
    @("to-webapp", $null),
    @("firewall2", $null),
    @("to-app", $null),
    @("webapp", $null),
    @("to-firewall2", $null),
    @("user-data", "X"),
    @("dmz", $null)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $a = $values[$i][0]
    $b = $values[$i][1]

    $ws.Cells.Item($row, 1).Value = $a
    if ($b -ne $null) {
        $ws.Cells.Item($row, 2).Value = $b
    } else {
        $ws.Cells.Item($row, 2).Value = $null
    }
}
